$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the "Driver Vintage" date-like strings stay as plain text (not auto-converted to dates)
$dateTextCells = @("E18","E19","E20","E21","E22","E23","E24","E25","E26","E27","E28","E29","E30","E31","E32","E33","E34","E35","E36","E37","E38","E39","E40","E41","E42","E43","E44","E45","E46","E47","E48","E49","E50","E51","E52","E53","E54","E55","E56","E57","E58","E59","E60","E61","E62","E63","E64")
foreach ($addr in $dateTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cell values from the weekly driver report refresh
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 12
$ws.Range("D3").Value = 78.59999999999999
$ws.Range("C4").Value = 29
$ws.Range("C5").Value = 91
$ws.Range("D6").Value = 98.7
$ws.Range("B7").Value = 22
$ws.Range("C7").Value = 190
$ws.Range("D7").Value = 98.8
$ws.Range("B8").Value = 45
$ws.Range("C8").Value = 328
$ws.Range("A16").Value = "Intel(R) Wireless-AC 9560 160MHz - 21.80.2.3"
$ws.Range("B16").Value = 10451
$ws.Range("D16").Value = 100
$ws.Range("A17").Value = "Intel(R) Wireless-AC 9560 160MHz - 21.10.2.2"
$ws.Range("B17").Value = 61902
$ws.Range("D17").Value = 100
$ws.Range("A18").Value = "Intel(R) Wireless-AC 9560 160MHz - 23.110.0.5"
$ws.Range("B18").Value = 54631
$ws.Range("D18").Value = 100
$ws.Range("E18").Value = "2025-01-01"
$ws.Range("A19").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B19").Value = 445055
$ws.Range("E19").Value = "2024-11-10"
$ws.Range("A20").Value = "Intel(R) Wireless-AC 9560 160MHz - 23.90.0.2"
$ws.Range("B20").Value = 4931894
$ws.Range("D20").Value = 99.90000000000001
$ws.Range("E20").Value = "2024-09-25"
$ws.Range("A21").Value = "Intel(R) Wireless-AC 9560 160MHz - 23.70.4.1"
$ws.Range("B21").Value = 52237
$ws.Range("E21").Value = "2024-08-13"
$ws.Range("A22").Value = "Intel(R) Wireless-AC 9560 160MHz - 23.50.0.6"
$ws.Range("B22").Value = 27295
$ws.Range("E22").Value = "2024-04-13"
$ws.Range("A23").Value = "Intel(R) Wireless-AC 9560 160MHz - 23.40.1.1"
$ws.Range("B23").Value = 276085
$ws.Range("D23").Value = 99.90000000000001
$ws.Range("E23").Value = "2024-03-19"
$ws.Range("A24").Value = "Intel(R) Wireless-AC 9560 160MHz - 23.30.0.6"
$ws.Range("B24").Value = 625139
$ws.Range("D24").Value = 99.90000000000001
$ws.Range("E24").Value = "2024-01-20"
$ws.Range("A25").Value = "Intel(R) Wireless-AC 9560 160MHz - 23.20.0.4"
$ws.Range("B25").Value = 44160
$ws.Range("D25").Value = 99.90000000000001
$ws.Range("E25").Value = "2023-11-28"
$ws.Range("A26").Value = "Intel(R) Wireless-AC 9560 160MHz - 23.10.0.8"
$ws.Range("B26").Value = 97122
$ws.Range("E26").Value = "2023-10-30"
$ws.Range("A27").Value = "Intel(R) Wireless-AC 9560 160MHz - 22.250.10.1"
$ws.Range("B27").Value = 78331
$ws.Range("E27").Value = "2023-08-14"
$ws.Range("A28").Value = "Intel(R) Wireless-AC 9560 160MHz - 22.220.0.4"
$ws.Range("B28").Value = 226852
$ws.Range("D28").Value = 99.90000000000001
$ws.Range("E28").Value = "2023-03-28"
$ws.Range("A29").Value = "Intel(R) Wireless-AC 9560 160MHz - 22.200.2.1"
$ws.Range("B29").Value = 453273
$ws.Range("E29").Value = "2023-03-08"
$ws.Range("A30").Value = "Intel(R) Wireless-AC 9560 160MHz - 22.190.0.4"
$ws.Range("B30").Value = 27599
$ws.Range("D30").Value = 100
$ws.Range("E30").Value = "2022-11-22"
$ws.Range("A31").Value = "Intel(R) Wireless-AC 9560 160MHz - 22.160.0.4"
$ws.Range("B31").Value = 3650830
$ws.Range("E31").Value = "2022-08-13"
$ws.Range("A32").Value = "Intel(R) Wireless-AC 9560 160MHz - 22.150.1.1"
$ws.Range("B32").Value = 154405
$ws.Range("E32").Value = "2022-06-20"
$ws.Range("A33").Value = "Intel(R) Wireless-AC 9560 160MHz - 22.100.0.3"
$ws.Range("B33").Value = 25808
$ws.Range("E33").Value = "2022-05-01"
$ws.Range("A34").Value = "Intel(R) Wireless-AC 9560 160MHz - 22.80.1.1"
$ws.Range("B34").Value = 94657
$ws.Range("D34").Value = 100
$ws.Range("E34").Value = "2022-05-01"
$ws.Range("A35").Value = "Intel(R) Wireless-AC 9560 160MHz - 22.140.0.3"
$ws.Range("B35").Value = 178916
$ws.Range("E35").Value = "2022-04-25"
$ws.Range("A36").Value = "Intel(R) Wireless-AC 9560 160MHz - 22.130.0.5"
$ws.Range("B36").Value = 109036
$ws.Range("D36").Value = 99.90000000000001
$ws.Range("E36").Value = "2022-03-14"
$ws.Range("A37").Value = "Intel(R) Wireless-AC 9560 160MHz - 22.120.0.3"
$ws.Range("B37").Value = 99547
$ws.Range("E37").Value = "2022-01-30"
$ws.Range("A38").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B38").Value = 77849
$ws.Range("E38").Value = "2021-08-18"
$ws.Range("A39").Value = "Intel(R) Wireless-AC 9560 160MHz - 22.70.0.6"
$ws.Range("B39").Value = 75637
$ws.Range("D39").Value = 100
$ws.Range("E39").Value = "2021-06-28"
$ws.Range("A40").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B40").Value = 34244
$ws.Range("D40").Value = 100
$ws.Range("E40").Value = "2021-04-27"
$ws.Range("A41").Value = "Intel(R) Wireless-AC 9560 160MHz - 22.50.0.7"
$ws.Range("B41").Value = 1543020
$ws.Range("E41").Value = "2021-04-18"
$ws.Range("A42").Value = "Intel(R) Wireless-AC 9560 160MHz - 22.40.0.7"
$ws.Range("B42").Value = 171222
$ws.Range("E42").Value = "2021-03-02"
$ws.Range("A43").Value = "Intel(R) Wireless-AC 9560 160MHz - 22.30.0.11"
$ws.Range("B43").Value = 238746
$ws.Range("D43").Value = 99.90000000000001
$ws.Range("E43").Value = "2021-01-19"
$ws.Range("A44").Value = "Intel(R) Wireless-AC 9560 160MHz - 22.10.0.7"
$ws.Range("B44").Value = 321983
$ws.Range("D44").Value = 99.90000000000001
$ws.Range("E44").Value = "2020-10-19"
$ws.Range("A45").Value = "Intel(R) Wireless-AC 9560 160MHz - 21.120.0.9"
$ws.Range("B45").Value = 95547
$ws.Range("E45").Value = "2020-08-15"
$ws.Range("A46").Value = "Intel(R) Wireless-AC 9560 160MHz - 21.110.3.2"
$ws.Range("B46").Value = 121232
$ws.Range("E46").Value = "2020-08-05"
$ws.Range("A47").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B47").Value = 59673
$ws.Range("D47").Value = 100
$ws.Range("E47").Value = "2020-08-05"
$ws.Range("A48").Value = "Intel(R) Wireless-AC 9560 160MHz - 21.110.2.1"
$ws.Range("B48").Value = 36791
$ws.Range("D48").Value = 99.90000000000001
$ws.Range("E48").Value = "2020-07-20"
$ws.Range("A49").Value = "Intel(R) Wireless-AC 9560 160MHz - 21.110.1.1"
$ws.Range("B49").Value = 684542
$ws.Range("D49").Value = 99.90000000000001
$ws.Range("E49").Value = "2020-06-30"
$ws.Range("A50").Value = "Intel(R) Wireless-AC 9560 160MHz - 21.90.1.2"
$ws.Range("B50").Value = 262523
$ws.Range("D50").Value = 100
$ws.Range("E50").Value = "2020-04-05"
$ws.Range("A51").Value = "Intel(R) Wireless-AC 9560 160MHz - 21.80.2.1"
$ws.Range("B51").Value = 209593
$ws.Range("D51").Value = 99.90000000000001
$ws.Range("E51").Value = "2020-02-24"
$ws.Range("A52").Value = "Intel(R) Wireless-AC 9560 160MHz - 21.70.0.6"
$ws.Range("B52").Value = 67365
$ws.Range("D52").Value = 100
$ws.Range("E52").Value = "2020-01-06"
$ws.Range("A53").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B53").Value = 113652
$ws.Range("E53").Value = "2020-01-06"
$ws.Range("A54").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B54").Value = 56018
$ws.Range("D54").Value = 100
$ws.Range("E54").Value = "2019-12-14"
$ws.Range("A55").Value = "Intel(R) Wireless-AC 9560 160MHz - 21.50.1.1"
$ws.Range("B55").Value = 308468
$ws.Range("D55").Value = 100
$ws.Range("E55").Value = "2019-10-05"
$ws.Range("A56").Value = "Intel(R) Wireless-AC 9560 160MHz - 21.40.2.2"
$ws.Range("B56").Value = 139827
$ws.Range("E56").Value = "2019-08-31"
$ws.Range("A57").Value = "Intel(R) Wireless-AC 9560 160MHz - 21.0.0.5"
$ws.Range("B57").Value = 194659
$ws.Range("E57").Value = "2019-08-31"
$ws.Range("A58").Value = "Intel(R) Wireless-AC 9560 160MHz - 21.40.1.4"
$ws.Range("B58").Value = 34718
$ws.Range("E58").Value = "2019-08-14"
$ws.Range("A59").Value = "Intel(R) Wireless-AC 9560 160MHz - 21.10.1.2"
$ws.Range("B59").Value = 546751
$ws.Range("E59").Value = "2019-08-10"
$ws.Range("A60").Value = "Intel(R) Wireless-AC 9560 160MHz - 21.30.3.2"
$ws.Range("B60").Value = 443012
$ws.Range("E60").Value = "2019-07-06"
$ws.Range("A61").Value = "Intel(R) Wireless-AC 9560 160MHz - 21.10.0.5"
$ws.Range("B61").Value = 108633
$ws.Range("E61").Value = "2019-04-06"
$ws.Range("A62").Value = "Intel(R) Wireless-AC 9560 160MHz - 21.0.1.1"
$ws.Range("B62").Value = 191877
$ws.Range("E62").Value = "2019-03-24"
$ws.Range("A63").Value = "Intel(R) Wireless-AC 9560 160MHz - 20.110.0.3"
$ws.Range("B63").Value = 61072
$ws.Range("D63").Value = 99.90000000000001
$ws.Range("E63").Value = "2018-11-27"
$ws.Range("A64").Value = "Intel(R) Wireless-AC 9560 160MHz - 20.100.0.4"
$ws.Range("B64").Value = 108823
$ws.Range("E64").Value = "2018-10-31"

Write-Host "Applied weekly driver report update for 2025-04-20"